$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column layout change -------------------------------------------------
# The columns "Pareto", "Unit No", "Sn Chassis", "Sn Engine" and
# "Production Year" (E:I) are replaced by two new columns "SN" and "User"
# (E:F). The remaining block of columns (PO No. .. Acquisition Value,
# previously J:O) shifts left to G:L, keeping its values and styles
# (Capitalized Date keeps its date-formatted style).

# Move J1:O2 -> G1:L2 (value + formatting) in one shot so the date style on
# the "Capitalized Date" column travels with its data.
$ws.Range("J1:O2").Copy($ws.Range("G1"))

# New columns E (SN) / F (User) headers; F2 already holds the value that
# belongs there ("KLX1034 (Opsional)") and is left untouched
$ws.Range("E1").Value = "SN"
$ws.Range("F1").Value = "User"
$ws.Range("E2").Value = "1057453 (Opsional)"

# Clear out the now-unused trailing columns M:O
$ws.Range("M1:O2").Clear()

# Column widths for the shifted columns (G:L), matching their new content
# (values chosen so the engine's internal pixel rounding lands as close as
# possible to the target widths of 23.109375/18.5546875/15.109375/
# 10.88671875/20.44140625/18.77734375 characters)
$ws.Columns("G").ColumnWidth = 22.3333333333333
$ws.Columns("H").ColumnWidth = 17.6666666666667
$ws.Columns("I").ColumnWidth = 14.3333333333333
$ws.Columns("J").ColumnWidth = 10
$ws.Columns("K").ColumnWidth = 19.6666666666667
$ws.Columns("L").ColumnWidth = 18

# Reset the view: drop the scrolled "topLeftCell" and move the active
# selection to F9
$ws.Range("F9").Select()
